# Planning V2 made and GDD reworked for V2
# Applies the content edits from the commit "planning V2 made and GDD reworked for V2"
# to the Rick Beniers planning workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix a typo / wording tweak in the "documentatie V1 inleveren" entry ---
$ws.Range("C35").Value = "project : documentatie V1 inleveren, trello bijwerken"

# --- Re-word the in-game winkel task (now phrased as "onderzoeken hoe/wat") ---
$ws.Range("C37").Value = "GameDev : in-game winkel onderzoeken hoe/wat"

# --- Fill in the previously empty "opgeleverd" notes for the V2 handover row ---
$ws.Range("F33").Value = "OOP : verslag afgemaakt en ingeleverd"
$ws.Range("F34").Value = "GameDev : planning V2 en GDD V2 ingeleverd "

# --- Replace the generic "activiteiten bepaald na TestPlan 2" placeholders ---
# with the real sprint-4 planning content (Object pooling / base of operations / brandstof)
$ws.Range("C53").Value = "GameDev : Object pooling(OPO) onderzoeken hoe/wat"
$ws.Range("C54").Value = "(Project : sprite baseOfOperations implementeren, OOP : verder werken aan OPO)"
$ws.Range("C55").Value = "GameDev : Object pooling(OPO) afmaken"
$ws.Range("C57").Value = "OOP : base of operations werkend maken met de rest van de game"
$ws.Range("C58").Value = "(GameDev : brandstof toevoegen als limiterende factor"
$ws.Range("C59").Value = "project : brandstof werkend maken met de rest van de game"

# --- "Versie 2 inleveren" -> "Versie 2 opleveren + documentatie" ---
$ws.Range("C63").Value = "Versie 2 opleveren + documentatie"

# --- Update the view so the window shows the sprint-4 area, selection on C66 ---
$ws.Activate()
$ws.Range("C66").Select()
